# 11.10.2020 MC Sales Details
# Insert a new price-list entry ("D74", RP 890, CP 960) above the existing
# "D92" row (row 28), pushing the Model/RP/CP/Lifting-discount columns
# (A:D) for that block down by one row. The second mini-table in columns
# F:I (the "V../Z.." Lifting list) is independent and must stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the A:D block for rows 28-41 down to rows 29-42, working from the
# bottom up so we never overwrite a value before it has been copied down.
for ($r = 41; $r -ge 28; $r--) {
    $destRow = $r + 1
    for ($col = 1; $col -le 4; $col++) {
        $ws.Cells.Item($destRow, $col).Value2 = $ws.Cells.Item($r, $col).Value2
    }
}

# New first entry of the block: model "D74", RP 890, CP 960, no discount.
$ws.Range("A28").Value2 = "D74"
$ws.Range("B28").Value2 = 890
$ws.Range("C28").Value2 = 960
$ws.Range("D28").ClearContents()

# Restore the active selection left by the editor.
$ws.Range("L46").Select() | Out-Null
